# The "Recipes" sheet had its columns reordered from
#   A=name, B=id, C=ingredients, D=steps
# to
#   A=id, B=ingredients, C=steps, D=name
# Re-write each row with the values in their new column positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Recipes")

# Header row
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "ingredients"
$ws.Range("C1").Value = "steps"
$ws.Range("D1").Value = "name"

# Row 2 - Pasta
$ws.Range("A2").Value = 12
$ws.Range("B2").Value = "Pasta, Sauce, Water, Salt"
$ws.Range("C2").Value = "Boil water`nPut Pasta in water`nWait for pasta to cook`nAdd Sauce and salt`nEnjoy`n"
$ws.Range("D2").Value = "Pasta"

# Row 3 - Cake
$ws.Range("A3").Value = 13
$ws.Range("B3").Value = "Flour, Eggs, Milk, Frosting"
$ws.Range("C3").Value = "Combine Flour, Eggs, and Milk in Mixer`nBake Cake in oven`nAdd Frosting`n"
$ws.Range("D3").Value = "Cake"

# Row 4 - Eggs
$ws.Range("A4").Value = 14
$ws.Range("B4").Value = "Eggs, Salt, Pepper, Cheese"
$ws.Range("C4").Value = "Whisk Eggs`nHeat up pan`nPut eggs in pan and season with salt and pepper`nAdd cheese and melt`n"
$ws.Range("D4").Value = "Eggs"
